$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1827586206896552
$ws.Range("C2").Value = 0.5724137931034483
$ws.Range("J2").Value = 0.01379310344827586
$ws.Range("P2").Value = 0.1344827586206896
$ws.Range("S2").Value = 0.09655172413793103
$ws.Range("B3").Value = 0.01734104046242774
$ws.Range("C3").Value = 0.02312138728323699
$ws.Range("J3").Value = 0.02890173410404624
$ws.Range("P3").Value = 0.7398843930635838
$ws.Range("S3").Value = 0.1907514450867052
$ws.Range("J4").Value = 0.07843137254901961
$ws.Range("P4").Value = 0.6862745098039216
$ws.Range("S4").Value = 0.2352941176470588
$ws.Range("B6").Value = 0.04484304932735426
$ws.Range("D6").Value = 0.01345291479820628
$ws.Range("F6").Value = 0.08520179372197309
$ws.Range("J6").Value = 0.2869955156950673
$ws.Range("O6").Value = 0.008968609865470852
$ws.Range("Q6").Value = 0.2197309417040359
$ws.Range("R6").Value = 0.04035874439461883
$ws.Range("S6").Value = 0.3004484304932735
$ws.Range("B7").Value = 0.0846774193548387
$ws.Range("D7").Value = 0.02016129032258064
$ws.Range("F7").Value = 0.04435483870967742
$ws.Range("J7").Value = 0.1653225806451613
$ws.Range("O7").Value = 0.02016129032258064
$ws.Range("Q7").Value = 0.1935483870967742
$ws.Range("R7").Value = 0.04838709677419355
$ws.Range("S7").Value = 0.4233870967741936
$ws.Range("B8").Value = 0.09973753280839895
$ws.Range("D8").Value = 0.01312335958005249
$ws.Range("F8").Value = 0.06561679790026247
$ws.Range("J8").Value = 0.1443569553805774
$ws.Range("O8").Value = 0.02362204724409449
$ws.Range("Q8").Value = 0.2125984251968504
$ws.Range("R8").Value = 0.05774278215223097
$ws.Range("S8").Value = 0.3832020997375328
$ws.Range("B9").Value = 0.0975609756097561
$ws.Range("D9").Value = 0.01829268292682927
$ws.Range("F9").Value = 0.0426829268292683
$ws.Range("J9").Value = 0.1341463414634146
$ws.Range("O9").Value = 0.01219512195121951
$ws.Range("Q9").Value = 0.1951219512195122
$ws.Range("R9").Value = 0.06097560975609756
$ws.Range("S9").Value = 0.4390243902439024
$ws.Range("B10").Value = 0.1118935837245696
$ws.Range("D10").Value = 0.0297339593114241
$ws.Range("E10").Value = 0.002347417840375587
$ws.Range("F10").Value = 0.06885758998435054
$ws.Range("J10").Value = 0.1369327073552426
$ws.Range("O10").Value = 0.02190923317683881
$ws.Range("Q10").Value = 0.2175273865414711
$ws.Range("R10").Value = 0.05633802816901409
$ws.Range("S10").Value = 0.3544600938967136
$ws.Range("G11").Value = 0.1421188630490956
$ws.Range("J11").Value = 0.1136950904392765
$ws.Range("K11").Value = 0.1937984496124031
$ws.Range("L11").Value = 0.5322997416020672
$ws.Range("S11").Value = 0.01808785529715762
$ws.Range("G12").Value = 0.8215962441314554
$ws.Range("J12").Value = 0.136150234741784
$ws.Range("K12").Value = 0.01408450704225352
$ws.Range("L12").Value = 0.009389671361502348
$ws.Range("S12").Value = 0.0187793427230047
$ws.Range("G13").Value = 0.6181818181818182
$ws.Range("J13").Value = 0.3272727272727273
$ws.Range("S13").Value = 0.05454545454545454
$ws.Range("F15").Value = 0.02164502164502164
$ws.Range("H15").Value = 0.1341991341991342
$ws.Range("I15").Value = 0.08658008658008658
$ws.Range("J15").Value = 0.3463203463203463
$ws.Range("K15").Value = 0.0735930735930736
$ws.Range("M15").Value = 0.01298701298701299
$ws.Range("N15").Value = 0.004329004329004329
$ws.Range("O15").Value = 0.06926406926406926
$ws.Range("S15").Value = 0.2510822510822511
$ws.Range("F16").Value = 0.02525252525252525
$ws.Range("H16").Value = 0.1515151515151515
$ws.Range("I16").Value = 0.09090909090909091
$ws.Range("J16").Value = 0.3434343434343434
$ws.Range("K16").Value = 0.1868686868686869
$ws.Range("M16").Value = 0.01515151515151515
$ws.Range("O16").Value = 0.0505050505050505
$ws.Range("S16").Value = 0.1363636363636364
$ws.Range("F17").Value = 0.0285132382892057
$ws.Range("H17").Value = 0.1771894093686354
$ws.Range("I17").Value = 0.06313645621181263
$ws.Range("J17").Value = 0.3686354378818738
$ws.Range("K17").Value = 0.1486761710794297
$ws.Range("M17").Value = 0.02443991853360489
$ws.Range("O17").Value = 0.07331975560081466
$ws.Range("S17").Value = 0.1160896130346232
$ws.Range("F18").Value = 0.008130081300813009
$ws.Range("H18").Value = 0.1382113821138211
$ws.Range("I18").Value = 0.08130081300813008
$ws.Range("J18").Value = 0.4390243902439024
$ws.Range("K18").Value = 0.1382113821138211
$ws.Range("M18").Value = 0.03252032520325204
$ws.Range("O18").Value = 0.04878048780487805
$ws.Range("S18").Value = 0.1138211382113821
$ws.Range("F19").Value = 0.01657000828500414
$ws.Range("H19").Value = 0.1822700911350456
$ws.Range("I19").Value = 0.07207953603976802
$ws.Range("J19").Value = 0.3728251864125932
$ws.Range("K19").Value = 0.1317315658657829
$ws.Range("M19").Value = 0.02734051367025684
$ws.Range("N19").Value = 0.001657000828500414
$ws.Range("O19").Value = 0.07705053852526926
$ws.Range("S19").Value = 0.1184755592377796
